# Lab04: update stack trace spreadsheet
# Adds a second (mirrored) stack-trace diagram in columns H:L (rows 77-99)
# documenting the square()/sumOfSquares()/main() call stack, alongside the
# existing diagram in columns A:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for the new diagram columns (I, J, L, M) ---------------
$ws.Range("I:I").ColumnWidth = 28.166666666666668
$ws.Range("J:J").ColumnWidth = 30.5
$ws.Range("L:L").ColumnWidth = 20.666666666666668
$ws.Range("M:M").ColumnWidth = 8.333333333333332

# --- "square(int x)" stack frame (rows 77-82) ------------------------------
$ws.Range("J77:J82").HorizontalAlignment = -4131
$ws.Range("J77:J82").Merge()
$ws.Range("J77").Value = "square(int x)"

$ws.Range("I78").Value = "int product"
$ws.Range("I79").Value = "Save R2"
$ws.Range("I80").Value = "Save R1"
$ws.Range("I81").Value = "previous frame pointer"
$ws.Range("K81").Value = "R5"
$ws.Range("L81").Value = "current frame pointer"
$ws.Range("I82").Value = "square() return address"
$ws.Range("K82").Value = "R7"
$ws.Range("L82").Value = "current return address"

# --- "sumOfSquares(int a[], int arraySize)" stack frame (rows 83-89) ------
$ws.Range("J83:J89").HorizontalAlignment = -4131
$ws.Range("J83:J89").Merge()
$ws.Range("J83").Value = "sumOfSquares(int a[], int arraySize)"

$ws.Range("I83").Value = "int x"
$ws.Range("I84").Value = "int sum"
$ws.Range("H85").Value = "x0000"
$ws.Range("I85").Value = "int counter"
$ws.Range("I86").Value = "Save R2"
$ws.Range("I87").Value = "Save R1"
$ws.Range("H88").Value = "x####"
$ws.Range("I88").Value = "previous frame pointer"
$ws.Range("K88").Value = "R5"
$ws.Range("L88").Value = "current frame pointer"
$ws.Range("I89").Value = "sumOfSquares() return address"
$ws.Range("K89").Value = "R7"
$ws.Range("L89").Value = "current return address"

# --- "main()" stack frame (rows 90-99) -------------------------------------
$ws.Range("J90:J99").HorizontalAlignment = -4131
$ws.Range("J90:J99").Merge()
$ws.Range("J90").Value = "main()"

$ws.Range("I90").Value = "array pointer param"
$ws.Range("H91").Value = "x0005"
$ws.Range("I91").Value = "arraySize param"
$ws.Range("H92").Value = "x0001"
$ws.Range("I92").Value = "array[0]"
$ws.Range("H93").Value = "x0000"
$ws.Range("I93").Value = "array[1]"
$ws.Range("H94").Value = "x0005"
$ws.Range("I94").Value = "array[2]"
$ws.Range("H95").Value = "x0003"
$ws.Range("I95").Value = "array[3]"
$ws.Range("H96").Value = "x0002"
$ws.Range("I96").Value = "array[4]"
$ws.Range("I97").Value = "int total"
$ws.Range("I98").Value = "previous frame pointer"
$ws.Range("K98").Value = "R5"
$ws.Range("L98").Value = "current frame pointer"
$ws.Range("I99").Value = "main() return address"
$ws.Range("K99").Value = "R7"
$ws.Range("L99").Value = "current return address"

# --- Sheet view: zoom + scroll position + selection ------------------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("N81").Select()

Write-Output "edit complete"
